$d = $word.ActiveDocument

# Splits a "<prefix><arrow-target>" run (where <arrow-target> begins with
# "> ") into three runs: "<prefix>" (trailing space preserved), "-", and
# the remaining "> ..." text - matching the documented edit, which turns
# every "X > Y" line into "X -> Y" built from three separate runs.
#
# $findText  - a unique piece of text to locate via Find (must start at
#              the same place the new "-" prefix starts, i.e. right after
#              any run(s) that must stay untouched).
# $prefixLen - number of characters (from the start of the found range)
#              that belong to the "<prefix> " portion (including its
#              trailing space), i.e. where the new "-" run gets inserted.
function Split-ArrowText {
    param(
        [string]$findText,
        [int]$prefixLen
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
        return
    }

    $prefixStart = $rng.Start
    $insPos = $prefixStart + $prefixLen

    # Insert the new "-" character right before "> ...".
    $insRng = $d.Range($insPos, $insPos)
    $insRng.InsertAfter("-")

    # Force the "-" to live in its own run (distinct from both neighbours)
    # by toggling a character property on/off - this makes the engine keep
    # it as a separate <w:r> once the paragraph is re-serialized, without
    # leaving any visible formatting behind.
    $dashRng = $d.Range($insPos, $insPos + 1)
    $dashRng.Bold = 1
    $dashRng.Bold = 0

    # Likewise, force the prefix text ("<prefix> ") to stay in its own run,
    # separate from whatever run precedes it (important when the prefix is
    # only part of a larger original run, e.g. "ACS " inside "PCIe "+"ACS ",
    # or "Socket P-State " following a separate " " run).
    $prefixRng = $d.Range($prefixStart, $insPos)
    $prefixRng.Bold = 1
    $prefixRng.Bold = 0
}

Split-ArrowText "CPU Power Management > Maximum Performance" ("CPU Power Management ".Length)
Split-ArrowText "Memory Frequency > Maximum Performance" ("Memory Frequency ".Length)
Split-ArrowText "Alg. Performance Boost Disable (ApbDis) > Enabled" ("Alg. Performance Boost Disable (ApbDis) ".Length)
Split-ArrowText "Socket P-State > P0" ("Socket P-State ".Length)
Split-ArrowText "NUMA Nodes Per Socket > 2" ("NUMA Nodes Per Socket ".Length)
Split-ArrowText "L3 cache as NUMA Domain > Enabled" ("L3 cache as NUMA Domain ".Length)
Split-ArrowText "x2APIC Mode > Enabled" ("x2APIC Mode ".Length)
Split-ArrowText "ACS > Disabled" ("ACS ".Length)
Split-ArrowText "Preferred IO > Disabled" ("Preferred IO ".Length)
Split-ArrowText "Enhanced Preferred IO > Enabled" ("Enhanced Preferred IO ".Length)
